# Semaine 5 update (v5.12): append two new bullet items to the list
# at the end of the document, after "Reprise d'OFv4, ...".
#
# Both new paragraphs re-use the existing list formatting
# (pStyle "Paragraphedeliste", numPr ilvl=0 / numId=1). The second one
# also carries a <w:lastRenderedPageBreak/> marker on its run, exactly
# as in the target OOXML, so we build the paragraphs from a raw
# WordprocessingML fragment and insert it at the very end of the
# document's main story (collapsed to a point, so InsertXML appends
# rather than replacing existing content).

$d = $word.ActiveDocument

$end = $d.Content.Duplicate
$end.Collapse(0)

$newItemsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Correction de la fonction remove de Kibana, reprise a 0 des ids, decallage des ids.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Correction d'asynchronisme get business fields dans la liste du cote resultats.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $end.InsertXML($newItemsXml)
